# Fix js calendar. Fix Bożena and 'Name, Surname' in resources
#
# 1) Normalize the professor name "Iwaniec Joanna" -> "Joanna Iwaniec "
#    (Surname-first -> "Name, Surname"-style ordering, with a trailing
#    space as authored).
# 2) De-duplicate "Giermek Bozena" (missing the Polish diacritic) into the
#    already-correct "Giermek Bożena" spelling used elsewhere in the sheet.
# 3) Leave the active selection on E48 (matching the saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Replace("Iwaniec Joanna", "Joanna Iwaniec ")
$ws.Cells.Replace("Giermek Bozena", "Giermek Bożena")

$ws.Range("E48").Select()
